{"js": "// Updates the division problems in the table cells, each cell's\n// text is an exact '<dividend>\u00f7<divisor>=' run; replace old -> new\n// on a one-to-one basis (matchCase so we don't clobber partial hits).\nconst replacements = [\n  [\"741\u00f73=\", \"236\u00f77=\"],\n  [\"245\u00f74=\", \"620\u00f75=\"],\n  [\"843\u00f73=\", \"658\u00f76=\"],\n  [\"723\u00f78=\", \"522\u00f74=\"],\n  [\"146\u00f73=\", \"205\u00f79=\"],\n  [\"566\u00f79=\", \"794\u00f74=\"],\n  [\"878\u00f73=\", \"772\u00f75=\"],\n  [\"998\u00f78=\", \"968\u00f72=\"],\n  [\"741\u00f75=\", \"197\u00f75=\"],\n  [\"570\u00f73=\", \"645\u00f75=\"],\n  [\"455\u00f73=\", \"762\u00f78=\"],\n  [\"679\u00f76=\", \"825\u00f78=\"],\n  [\"242\u00f73=\", \"666\u00f74=\"],\n  [\"377\u00f77=\", \"693\u00f77=\"],\n  [\"882\u00f72=\", \"215\u00f77=\"],\n  [\"110\u00f79=\", \"445\u00f73=\"],\n  [\"481\u00f73=\", \"672\u00f78=\"],\n  [\"203\u00f75=\", \"700\u00f74=\"],\n  [\"448\u00f74=\", \"444\u00f74=\"],\n  [\"960\u00f75=\", \"424\u00f73=\"],\n  [\"771\u00f77=\", \"697\u00f79=\"],\n  [\"767\u00f72=\", \"480\u00f78=\"],\n  [\"984\u00f74=\", \"906\u00f73=\"],\n  [\"474\u00f74=\", \"406\u00f78=\"],\n  [\"342\u00f78=\", \"142\u00f73=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Search text not found: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n", "ps1": "# Updates the division problems in the table cells. Each target cell\n# holds an exact '<dividend>\u00f7<divisor>=' run; replace old -> new using\n# Find/Replace across the whole document body, one pair at a time so we\n# never touch an already-updated cell.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"741\u00f73=\", \"236\u00f77=\"),\n    @(\"245\u00f74=\", \"620\u00f75=\"),\n    @(\"843\u00f73=\", \"658\u00f76=\"),\n    @(\"723\u00f78=\", \"522\u00f74=\"),\n    @(\"146\u00f73=\", \"205\u00f79=\"),\n    @(\"566\u00f79=\", \"794\u00f74=\"),\n    @(\"878\u00f73=\", \"772\u00f75=\"),\n    @(\"998\u00f78=\", \"968\u00f72=\"),\n    @(\"741\u00f75=\", \"197\u00f75=\"),\n    @(\"570\u00f73=\", \"645\u00f75=\"),\n    @(\"455\u00f73=\", \"762\u00f78=\"),\n    @(\"679\u00f76=\", \"825\u00f78=\"),\n    @(\"242\u00f73=\", \"666\u00f74=\"),\n    @(\"377\u00f77=\", \"693\u00f77=\"),\n    @(\"882\u00f72=\", \"215\u00f77=\"),\n    @(\"110\u00f79=\", \"445\u00f73=\"),\n    @(\"481\u00f73=\", \"672\u00f78=\"),\n    @(\"203\u00f75=\", \"700\u00f74=\"),\n    @(\"448\u00f74=\", \"444\u00f74=\"),\n    @(\"960\u00f75=\", \"424\u00f73=\"),\n    @(\"771\u00f77=\", \"697\u00f79=\"),\n    @(\"767\u00f72=\", \"480\u00f78=\"),\n    @(\"984\u00f74=\", \"906\u00f73=\"),\n    @(\"474\u00f74=\", \"406\u00f78=\"),\n    @(\"342\u00f78=\", \"142\u00f73=\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1  # wdFindContinue\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($find.Text, $false, $true, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n\n"}
